$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2808873333333333
$ws.Range("H2").Value = 0.842662
$ws.Range("I2").Value = 0.5595554696739399
$ws.Range("J2").Value = 0.5595554696739399
$ws.Range("M2").Value = 1.979605666666667
$ws.Range("N2").Value = 5.938817
$ws.Range("O2").Value = 0.05865520690928468
$ws.Range("P2").Value = 0.05865520690928468
$ws.Range("Q2").Value = 0.5560461567615556
$ws.Range("R2").Value = 5.004415410854
$ws.Range("S2").Value = 0.03282084185094691
$ws.Range("T2").Value = 0.03282084185094691

$ws.Range("G3").Value = 0.2808873333333333
$ws.Range("H3").Value = 0.842662
$ws.Range("I3").Value = 0.5595554696739399
$ws.Range("J3").Value = 0.5595554696739399
$ws.Range("M3").Value = 16.37791666666667
$ws.Range("N3").Value = 49.13375000000001
$ws.Range("O3").Value = 0.4852734597612733
$ws.Range("P3").Value = 0.4852734597612734
$ws.Range("Q3").Value = 4.600349338055556
$ws.Range("R3").Value = 41.40314404250001
$ws.Range("S3").Value = 0.271537418697017
$ws.Range("T3").Value = 0.2715374186970171

$ws.Range("G4").Value = 0.2808873333333333
$ws.Range("H4").Value = 0.842662
$ws.Range("I4").Value = 0.5595554696739399
$ws.Range("J4").Value = 0.5595554696739399
$ws.Range("M4").Value = 15.39234866666667
$ws.Range("N4").Value = 46.177046
$ws.Range("O4").Value = 0.4560713333294419
$ws.Range("P4").Value = 0.4560713333294419
$ws.Range("Q4").Value = 4.323515770716889
$ws.Range("R4").Value = 38.91164193645201
$ws.Range("S4").Value = 0.2551972091259759
$ws.Range("T4").Value = 0.2551972091259759

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2210956666666667
$ws.Range("H5").Value = 0.663287
$ws.Range("I5").Value = 0.4404445303260602
$ws.Range("J5").Value = 0.4404445303260602
$ws.Range("M5").Value = 1.979605666666667
$ws.Range("N5").Value = 5.938817
$ws.Range("O5").Value = 0.05865520690928468
$ws.Range("P5").Value = 0.05865520690928468
$ws.Range("Q5").Value = 0.4376822346087778
$ws.Range("R5").Value = 3.939140111479
$ws.Range("S5").Value = 0.02583436505833777
$ws.Range("T5").Value = 0.02583436505833777

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2210956666666667
$ws.Range("H6").Value = 0.663287
$ws.Range("I6").Value = 0.4404445303260602
$ws.Range("J6").Value = 0.4404445303260602
$ws.Range("M6").Value = 16.37791666666667
$ws.Range("N6").Value = 49.13375000000001
$ws.Range("O6").Value = 0.4852734597612733
$ws.Range("P6").Value = 0.4852734597612734
$ws.Range("Q6").Value = 3.621086404027778
$ws.Range("R6").Value = 32.58977763625001
$ws.Range("S6").Value = 0.2137360410642563
$ws.Range("T6").Value = 0.2137360410642563

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2210956666666667
$ws.Range("H7").Value = 0.663287
$ws.Range("I7").Value = 0.4404445303260602
$ws.Range("J7").Value = 0.4404445303260602
$ws.Range("M7").Value = 15.39234866666667
$ws.Range("N7").Value = 46.177046
$ws.Range("O7").Value = 0.4560713333294419
$ws.Range("P7").Value = 0.4560713333294419
$ws.Range("Q7").Value = 3.403181590022445
$ws.Range("R7").Value = 30.628634310202
$ws.Range("S7").Value = 0.2008741242034661
$ws.Range("T7").Value = 0.2008741242034661
